$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.260.13'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '1.650.79'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.90'
$ws.Range('E5').Value = '  -0.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.511'
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.08'
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').Value = '1.885.65'
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('D13').Value = '1.658.44'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('E14').Value = '  -1.36%  '
$ws.Range('E15').Value = '  +2.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.73'
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('D17').Value = '27.283.46'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '220.22'
$ws.Range('E19').Value = '  -1.66%  '
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.83'
$ws.Range('E21').Value = '  +1.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.55'
$ws.Range('E22').Value = '  +5.57%  '
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.21'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.10'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.56'
$ws.Range('E26').Value = '  +1.53%  '
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  -0.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.83'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0510'
$ws.Range('E30').Value = '  -1.18%  '
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('E32').Value = '  -1.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.03'
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('E34').Value = '  +0.82%  '
$ws.Range('D35').Value = '1.259.00'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.547'
$ws.Range('E38').Value = '  +1.65%  '
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.45'
$ws.Range('E42').Value = '  +1.26%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.23'
$ws.Range('E43').Value = '  +5.04%  '
$ws.Range('D44').Value = '1.795.32'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.09'
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.74'
$ws.Range('E46').Value = '  -0.20%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('E48').Value = '  +22.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0513'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0973'
$ws.Range('E51').Value = '  -0.94%  '
